# GeoMap in proxidrugs_py updated
#
# Insert a new "location" column (clean city name, no "Städte" suffix)
# right after the Institute/Company column, populate it for the Hamburg
# rows, and add a new "counts" column that excludes the Hamburg / München
# rows, finishing with a SUM total under it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (old B "Location" -> C, old C "counts" -> D).
# Excel's column insert duplicates the left neighbour's formatting, so
# match column A's width on the new column B.
$ws.Columns("B:B").Insert()
$ws.Columns("B:B").ColumnWidth = 29

# New column B header + values (only populated for the Hamburg Städte rows)
$ws.Range("B1").Value = "location"
$ws.Range("B3").Value = "Hamburg"
$ws.Range("B9").Value = "Hamburg"

# New column E: copy of the counts (col D), skipping the Hamburg Städte
# (rows 3 & 9) and München (row 6) rows, with a sum total below.
$ws.Range("E2").Value = 7
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 35
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 10
$ws.Range("E10").Value = 6
$ws.Range("E11").Formula = "=SUM(E2:E10)"

[void]$ws.Range("A10").Select()
